$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data-wide-value")
$ws.Columns.Item(2).Delete()
